# Update countries & provincias Spain
# Refreshes the "Pais" sheet with newer source data: the table was re-sorted
# (Bolivia and Venezuela moved up to reflect their new, higher case counts)
# and several countries' case figures were updated, plus the "last updated"
# timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 03:22"

# --- Panama (row 49) updated case counts -------------------------------
$ws.Range("B49").Value = 5166
$ws.Range("C49").Value = 174
$ws.Range("D49").Value = 271
$ws.Range("E49").Value = 4749
$ws.Range("F49").Value = 86
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 146

# --- Bolivia moves up (now row 95) with refreshed figures --------------
$ws.Range("A95").Value = "Bolivia"
$ws.Range("B95").Value = 703
$ws.Range("C95").Value = 31
$ws.Range("D95").Value = 44
$ws.Range("E95").Value = 616
$ws.Range("F95").Value = 3
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 43

# --- Libano shifts down to row 96 (figures unchanged) -------------------
$ws.Range("A96").Value = "Libano"
$ws.Range("B96").Value = 688
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 140
$ws.Range("E96").Value = 526
$ws.Range("F96").Value = 46
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 22

# --- Costa Rica shifts down to row 97 (figures unchanged) ---------------
$ws.Range("A97").Value = "Costa Rica"
$ws.Range("B97").Value = 686
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 196
$ws.Range("E97").Value = 484
$ws.Range("F97").Value = 8
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 6

# --- Uruguay (row 102) updated case counts ------------------------------
$ws.Range("B102").Value = 557
$ws.Range("C102").Value = 8
$ws.Range("D102").Value = 354
$ws.Range("E102").Value = 191

# --- Venezuela moves up (now row 120) with refreshed figures ------------
$ws.Range("A120").Value = "Venezuela"
$ws.Range("B120").Value = 311
$ws.Range("C120").Value = 13
$ws.Range("D120").Value = 126
$ws.Range("E120").Value = 175
$ws.Range("F120").Value = 4
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 10

# --- Mali shifts down to row 121 (figures unchanged) ---------------------
$ws.Range("A121").Value = "Mali"
$ws.Range("B121").Value = 309
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 77
$ws.Range("E121").Value = 211
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 21

# --- Isla de Man shifts down to row 122 (figures unchanged) --------------
$ws.Range("A122").Value = "Isla de Man"
$ws.Range("B122").Value = 307
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 221
$ws.Range("E122").Value = 70
$ws.Range("F122").Value = 20
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 16
